# "Completed all fake data inserts"
#
# The "Inserts >>" sheet was still empty. Finish it off the same way the
# other helper sheets ("Store Procs", "Rollback", ...) are built: a header
# row (FileName / Schema / Type / Command) followed by one sample/"fake"
# data row whose Command column builds the sqlcmd invocation from the
# named ranges defined on the "Entries" sheet.

$wb = $excel.ActiveWorkbook

$dst = $wb.Worksheets.Item("Inserts >>")

# Header row
$dst.Range("A1").Value = "FileName"
$dst.Range("B1").Value = "Schema"
$dst.Range("C1").Value = "Type"
$dst.Range("D1").Value = "Command"

# Sample/fake data row
$dst.Range("A2").Value = "SP_DeleteEmails"
$dst.Range("B2").Value = "cust."
$dst.Range("C2").Value = ".StoredProcedure.sql"
$dst.Range("D2").Formula = '="sqlcmd -S "&Entries_Server&" -d "&Entries_Database&" -U "&Entries_User&" -P "&Entries_Password&" -i "&Entries_StoredProcedure_RootPath&B2&A2&C2&""""&" >> """&Entries_StoredProcedure_DebugText'

# Leave the same A1:D2 / A1:C2 ranges selected on the sheets this was
# copied/checked against while putting the data together.
$storeProcs = $wb.Worksheets.Item("Store Procs")
$storeProcs.Range("A1:D2").Select()

$rollback = $wb.Worksheets.Item("Rollback")
$rollback.Range("A1:C2").Select()

# Finish with "Inserts >>" as the active tab/selection.
$dst.Activate()
$dst.Range("A1:D2").Select()
